# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a number of rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 11;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 13;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 16;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 24;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 52;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 59;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 61;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 63;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 70;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 80;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 84;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 105; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 126; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 128; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 144; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 167; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
